$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("D2").Value = 44224
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 16500
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 16750
$ws.Range("S2").Value = 1047

# --- Row 3 updates ---
$ws.Range("D3").Value = 44224
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 14500
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14750
$ws.Range("S3").Value = 922

# --- Row 4 updates ---
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 12500
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 12750
$ws.Range("S4").Value = 797

# --- Row 5 updates ---
$ws.Range("D5").Value = 44209
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 15500
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 15750
$ws.Range("S5").Value = 984

# --- Row 6 updates ---
$ws.Range("D6").Value = 44594
$ws.Range("L6").Value = "Especial"
$ws.Range("M6").Value = 240
$ws.Range("N6").Value = 15500
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15750
$ws.Range("Q6").Value = "$/caja 15 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1050
$ws.Range("T6").Value = 15

# --- Row 7 updates ---
$ws.Range("D7").Value = 44594
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 13500
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 13750
$ws.Range("Q7").Value = "$/caja 15 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 917
$ws.Range("T7").Value = 15

# --- Row 8 updates ---
$ws.Range("L8").Value = "Primera"
$ws.Range("N8").Value = 16500
$ws.Range("O8").Value = 17000
$ws.Range("P8").Value = 16750
$ws.Range("S8").Value = 931

# --- Row 9 updates ---
$ws.Range("D9").Value = 44230
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 160
$ws.Range("N9").Value = 14500
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14750
$ws.Range("Q9").Value = "$/caja 18 kilos granel"
$ws.Range("S9").Value = 819
$ws.Range("T9").Value = 18

# --- New row 10 ---
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("D10").Value = 44210
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100103
$ws.Range("H10").Value = "Frutos de hueso (carozo)"
$ws.Range("I10").Value = 100103002
$ws.Range("J10").Value = "Ciruela"
$ws.Range("K10").Value = "Black Amber"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 240
$ws.Range("N10").Value = 15500
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15750
$ws.Range("Q10").Value = "$/caja 16 kilos granel"
$ws.Range("R10").Value = "Región Metropolitana"
$ws.Range("S10").Value = 984
$ws.Range("T10").Value = 16

# --- New row 11 ---
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("D11").Value = 44210
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100103
$ws.Range("H11").Value = "Frutos de hueso (carozo)"
$ws.Range("I11").Value = 100103002
$ws.Range("J11").Value = "Ciruela"
$ws.Range("K11").Value = "Black Amber"
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 12500
$ws.Range("O11").Value = 13000
$ws.Range("P11").Value = 12750
$ws.Range("Q11").Value = "$/caja 16 kilos granel"
$ws.Range("R11").Value = "Región Metropolitana"
$ws.Range("S11").Value = 797
$ws.Range("T11").Value = 16
